$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.504.75'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.905.78'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.46'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4847'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4074'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08149'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.49'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +5.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.915.75'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.030'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.081'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.56%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.40'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.43%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06745'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001044'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.72'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.511.54'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.581'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.96%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.83'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.67%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.080.52'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.19'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.08'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.284'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +10.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.104'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.06'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.039'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09555'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.530'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.396'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.54%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02266'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06113'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.172'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5978'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.919'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.46'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.57%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1857'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.426'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.280'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.07725'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5571'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.966'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '115.01'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.70'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.054'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.55%  '
